$wb = $excel.ActiveWorkbook

# "Ready for handoff" status values are now reported as "In Translation"
# wherever they occur (Overview + each per-language sheet).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Put the literal on the left: PowerShell's `-eq` coerces the
        # right-hand side to the left operand's type, and boolean cells
        # ("True"/"False" text) would otherwise coerce any non-empty
        # string to $true and false-match here.
        if ("Ready for handoff" -eq $cell.Text) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the per-language status columns to match the new report layout.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
